$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.080076666666666
$ws.Range("H2").Value = 6.240229999999999
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.848466999999999
$ws.Range("N2").Value = 29.545401
$ws.Range("O2").Value = 0.0466568297496787
$ws.Range("P2").Value = 0.0466568297496787
$ws.Range("Q2").Value = 20.48556640913666
$ws.Range("R2").Value = 184.37009768223
$ws.Range("S2").Value = 0.0466568297496787
$ws.Range("T2").Value = 0.0466568297496787

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.080076666666666
$ws.Range("H3").Value = 6.240229999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 198.2465873333333
$ws.Range("N3").Value = 594.7397619999999
$ws.Range("O3").Value = 0.9391875175767094
$ws.Range("P3").Value = 0.9391875175767094
$ws.Range("Q3").Value = 412.3681005583621
$ws.Range("R3").Value = 3711.312905025259
$ws.Range("S3").Value = 0.9391875175767094
$ws.Range("T3").Value = 0.9391875175767094

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.080076666666666
$ws.Range("H4").Value = 6.240229999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.988018666666667
$ws.Range("N4").Value = 8.964055999999999
$ws.Range("O4").Value = 0.01415565267361191
$ws.Range("P4").Value = 0.01415565267361191
$ws.Range("Q4").Value = 6.215307908097777
$ws.Range("R4").Value = 55.93777117287999
$ws.Range("S4").Value = 0.01415565267361191
$ws.Range("T4").Value = 0.01415565267361191
